$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last 3 data rows (previously rows 8-10, for the same-cluster
# "Sending cluster == Target cluster" pairs) so the data now only contains
# the 6 off-diagonal Sending/Target cluster combinations, re-computed with
# the new TPM values.
$ws.Rows.Item(8).Resize(3).Delete() | Out-Null

# New data (rows 2-7), columns A through T, reflecting the updated TPM
# calculation.
$data = @(
    @("ECs", "Ccl25", "Ccr10", "FAPs", 3, 1, 2.758234, 8.274702, 0.1921318935545868, 0.1921318935545868, 2, 0.6666666666666666, 1.238415, 3.715245, 0.4359607654144799, 0.4359607654144798, 3.41583835911, 30.74254523199, 0.08376196737459105, 0.08376196737459103),
    @("ECs", "Ccl25", "Ccr10", "MuSCs", 3, 1, 2.758234, 8.274702, 0.1921318935545868, 0.1921318935545868, 3, 1, 1.602242, 4.806725999999999, 0.5640392345855201, 0.5640392345855201, 4.419358360627999, 39.77422524565199, 0.1083699261799958, 0.1083699261799958),
    @("FAPs", "Ccl25", "Ccr10", "FAPs", 3, 1, 7.720664, 23.161992, 0.5378027367579149, 0.5378027367579149, 2, 0.6666666666666666, 1.238415, 3.715245, 0.4359607654144799, 0.4359607654144798, 9.561386107560001, 86.05247496804002, 0.2344608927589826, 0.2344608927589826),
    @("FAPs", "Ccl25", "Ccr10", "MuSCs", 3, 1, 7.720664, 23.161992, 0.5378027367579149, 0.5378027367579149, 3, 1, 1.602242, 4.806725999999999, 0.5640392345855201, 0.5640392345855201, 12.370372128688, 111.333349158192, 0.3033418439989323, 0.3033418439989323),
    @("MuSCs", "Ccl25", "Ccr10", "FAPs", 3, 1, 3.877042333333333, 11.631127, 0.2700653696874982, 0.2700653696874982, 2, 0.6666666666666666, 1.238415, 3.715245, 0.4359607654144799, 0.4359607654144798, 4.801387381235, 43.212486431115, 0.1177379052809062, 0.1177379052809062),
    @("MuSCs", "Ccl25", "Ccr10", "MuSCs", 3, 1, 3.877042333333333, 11.631127, 0.2700653696874982, 0.2700653696874982, 3, 1, 1.602242, 4.806725999999999, 0.5640392345855201, 0.5640392345855201, 6.211960062244666, 55.90764056020199, 0.152327464406592, 0.152327464406592)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
